# Weekly refresh: insert the newest week's two records (Primera / Segunda
# quality grades) at the top of the historical block for this market +
# variety, pushing the rest of the history down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 180, shifting existing rows 180-292
# down to 182-294 (and carrying the row-180 formatting, e.g. the date
# style on column D, onto the freshly inserted rows).
$ws.Range("A180:A181").EntireRow.Insert()

# Populate the new row 180 (Primera).
$ws.Cells.Item(180, 1).Value = 7
$ws.Cells.Item(180, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(180, 3).Value = "Ñuble"
$ws.Cells.Item(180, 4).Value = 44879
$ws.Cells.Item(180, 5).Value = 16
$ws.Cells.Item(180, 6).Value = 100112006
$ws.Cells.Item(180, 7).Value = "Repollo"
$ws.Cells.Item(180, 8).Value = "Crespo record"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 400
$ws.Cells.Item(180, 11).Value = 1500
$ws.Cells.Item(180, 12).Value = 1600
$ws.Cells.Item(180, 13).Value = 1550
$ws.Cells.Item(180, 14).Value = "$/unidad"
$ws.Cells.Item(180, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(180, 16).Value = 1550
$ws.Cells.Item(180, 17).Value = 1
$ws.Cells.Item(180, 18).Value = "Hortaliza"

# Populate the new row 181 (Segunda).
$ws.Cells.Item(181, 1).Value = 7
$ws.Cells.Item(181, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(181, 3).Value = "Ñuble"
$ws.Cells.Item(181, 4).Value = 44879
$ws.Cells.Item(181, 5).Value = 16
$ws.Cells.Item(181, 6).Value = 100112006
$ws.Cells.Item(181, 7).Value = "Repollo"
$ws.Cells.Item(181, 8).Value = "Crespo record"
$ws.Cells.Item(181, 9).Value = "Segunda"
$ws.Cells.Item(181, 10).Value = 300
$ws.Cells.Item(181, 11).Value = 1200
$ws.Cells.Item(181, 12).Value = 1300
$ws.Cells.Item(181, 13).Value = 1250
$ws.Cells.Item(181, 14).Value = "$/unidad"
$ws.Cells.Item(181, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(181, 16).Value = 1250
$ws.Cells.Item(181, 17).Value = 1
$ws.Cells.Item(181, 18).Value = "Hortaliza"
